$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 6 de Agosto de 2020 a las 20:01"

# Update country rows: a handful of countries were re-inserted/re-ordered in
# the list (so the row previously holding one country now holds another),
# and every row carries refreshed COVID-19 statistics for the newer data pull.

$ws.Cells.Item(4, 1).Value = "Estados Unidos"
$ws.Cells.Item(4, 2).Value = 4997868
$ws.Cells.Item(4, 3).Value = 24300
$ws.Cells.Item(4, 4).Value = 2549791
$ws.Cells.Item(4, 5).Value = 2286012
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 464
$ws.Cells.Item(4, 8).Value = 162065

$ws.Cells.Item(6, 1).Value = "India"
$ws.Cells.Item(6, 2).Value = 2025338
$ws.Cells.Item(6, 3).Value = 62099
$ws.Cells.Item(6, 4).Value = 1377384
$ws.Cells.Item(6, 5).Value = 606320
$ws.Cells.Item(6, 6).Value = 0
$ws.Cells.Item(6, 7).Value = 895
$ws.Cells.Item(6, 8).Value = 41634

$ws.Cells.Item(20, 1).Value = "Turquia"
$ws.Cells.Item(20, 2).Value = 237265
$ws.Cells.Item(20, 3).Value = 1153
$ws.Cells.Item(20, 4).Value = 220546
$ws.Cells.Item(20, 5).Value = 10921
$ws.Cells.Item(20, 6).Value = 0
$ws.Cells.Item(20, 7).Value = 14
$ws.Cells.Item(20, 8).Value = 5798

$ws.Cells.Item(22, 1).Value = "Alemania"
$ws.Cells.Item(22, 2).Value = 215100
$ws.Cells.Item(22, 3).Value = 996
$ws.Cells.Item(22, 4).Value = 196200
$ws.Cells.Item(22, 5).Value = 9651
$ws.Cells.Item(22, 6).Value = 0
$ws.Cells.Item(22, 7).Value = 4
$ws.Cells.Item(22, 8).Value = 9249

$ws.Cells.Item(27, 1).Value = "Canada"
$ws.Cells.Item(27, 2).Value = 118417
$ws.Cells.Item(27, 3).Value = 230
$ws.Cells.Item(27, 4).Value = 102947
$ws.Cells.Item(27, 5).Value = 6507
$ws.Cells.Item(27, 6).Value = 0
$ws.Cells.Item(27, 7).Value = 1
$ws.Cells.Item(27, 8).Value = 8963

$ws.Cells.Item(36, 1).Value = "Israel"
$ws.Cells.Item(36, 2).Value = 79275
$ws.Cells.Item(36, 3).Value = 1356
$ws.Cells.Item(36, 4).Value = 53412
$ws.Cells.Item(36, 5).Value = 25287
$ws.Cells.Item(36, 6).Value = 0
$ws.Cells.Item(36, 7).Value = 11
$ws.Cells.Item(36, 8).Value = 576

$ws.Cells.Item(59, 1).Value = "Argelia"
$ws.Cells.Item(59, 2).Value = 33626
$ws.Cells.Item(59, 3).Value = 571
$ws.Cells.Item(59, 4).Value = 23238
$ws.Cells.Item(59, 5).Value = 9115
$ws.Cells.Item(59, 6).Value = 0
$ws.Cells.Item(59, 7).Value = 12
$ws.Cells.Item(59, 8).Value = 1273

$ws.Cells.Item(60, 1).Value = "Azerbaiyan"
$ws.Cells.Item(60, 2).Value = 33247
$ws.Cells.Item(60, 3).Value = 144
$ws.Cells.Item(60, 4).Value = 29275
$ws.Cells.Item(60, 5).Value = 3493
$ws.Cells.Item(60, 6).Value = 0
$ws.Cells.Item(60, 7).Value = 3
$ws.Cells.Item(60, 8).Value = 479

$ws.Cells.Item(61, 1).Value = "Marruecos"
$ws.Cells.Item(61, 2).Value = 29644
$ws.Cells.Item(61, 3).Value = 1144
$ws.Cells.Item(61, 4).Value = 20553
$ws.Cells.Item(61, 5).Value = 8642
$ws.Cells.Item(61, 6).Value = 0
$ws.Cells.Item(61, 7).Value = 14
$ws.Cells.Item(61, 8).Value = 449

$ws.Cells.Item(66, 1).Value = "Kenia"
$ws.Cells.Item(66, 2).Value = 24411
$ws.Cells.Item(66, 3).Value = 538
$ws.Cells.Item(66, 4).Value = 10444
$ws.Cells.Item(66, 5).Value = 13568
$ws.Cells.Item(66, 6).Value = 0
$ws.Cells.Item(66, 7).Value = 8
$ws.Cells.Item(66, 8).Value = 399

$ws.Cells.Item(70, 1).Value = "Etiopia"
$ws.Cells.Item(70, 2).Value = 20900
$ws.Cells.Item(70, 3).Value = 564
$ws.Cells.Item(70, 4).Value = 9027
$ws.Cells.Item(70, 5).Value = 11508
$ws.Cells.Item(70, 6).Value = 0
$ws.Cells.Item(70, 7).Value = 9
$ws.Cells.Item(70, 8).Value = 365

$ws.Cells.Item(71, 1).Value = "Costa Rica"
$ws.Cells.Item(71, 2).Value = 20417
$ws.Cells.Item(71, 3).Value = 0
$ws.Cells.Item(71, 4).Value = 6851
$ws.Cells.Item(71, 5).Value = 13375
$ws.Cells.Item(71, 6).Value = 0
$ws.Cells.Item(71, 7).Value = 0
$ws.Cells.Item(71, 8).Value = 191

$ws.Cells.Item(79, 1).Value = "Estado de Palestina"
$ws.Cells.Item(79, 2).Value = 13398
$ws.Cells.Item(79, 3).Value = 333
$ws.Cells.Item(79, 4).Value = 6907
$ws.Cells.Item(79, 5).Value = 6400
$ws.Cells.Item(79, 6).Value = 0
$ws.Cells.Item(79, 7).Value = 2
$ws.Cells.Item(79, 8).Value = 91

$ws.Cells.Item(95, 1).Value = "Zambia"
$ws.Cells.Item(95, 2).Value = 7164
$ws.Cells.Item(95, 3).Value = 142
$ws.Cells.Item(95, 4).Value = 5786
$ws.Cells.Item(95, 5).Value = 1179
$ws.Cells.Item(95, 6).Value = 0
$ws.Cells.Item(95, 7).Value = 23
$ws.Cells.Item(95, 8).Value = 199

$ws.Cells.Item(96, 1).Value = "Luxemburgo"
$ws.Cells.Item(96, 2).Value = 7073
$ws.Cells.Item(96, 3).Value = 66
$ws.Cells.Item(96, 4).Value = 5750
$ws.Cells.Item(96, 5).Value = 1204
$ws.Cells.Item(96, 6).Value = 0
$ws.Cells.Item(96, 7).Value = 1
$ws.Cells.Item(96, 8).Value = 119

$ws.Cells.Item(105, 1).Value = "Maldivas"
$ws.Cells.Item(105, 2).Value = 4680
$ws.Cells.Item(105, 3).Value = 86
$ws.Cells.Item(105, 4).Value = 2725
$ws.Cells.Item(105, 5).Value = 1936
$ws.Cells.Item(105, 6).Value = 0
$ws.Cells.Item(105, 7).Value = 0
$ws.Cells.Item(105, 8).Value = 19

$ws.Cells.Item(106, 1).Value = "Republica de Africa Central"
$ws.Cells.Item(106, 2).Value = 4620
$ws.Cells.Item(106, 3).Value = 2
$ws.Cells.Item(106, 4).Value = 1641
$ws.Cells.Item(106, 5).Value = 2920
$ws.Cells.Item(106, 6).Value = 0
$ws.Cells.Item(106, 7).Value = 0
$ws.Cells.Item(106, 8).Value = 59

$ws.Cells.Item(107, 1).Value = "Hungria"
$ws.Cells.Item(107, 2).Value = 4597
$ws.Cells.Item(107, 3).Value = 33
$ws.Cells.Item(107, 4).Value = 3463
$ws.Cells.Item(107, 5).Value = 534
$ws.Cells.Item(107, 6).Value = 0
$ws.Cells.Item(107, 7).Value = 1
$ws.Cells.Item(107, 8).Value = 600

$ws.Cells.Item(117, 1).Value = "Mayotte"
$ws.Cells.Item(117, 2).Value = 3042
$ws.Cells.Item(117, 3).Value = 11
$ws.Cells.Item(117, 4).Value = 2738
$ws.Cells.Item(117, 5).Value = 265
$ws.Cells.Item(117, 6).Value = 0
$ws.Cells.Item(117, 7).Value = 0
$ws.Cells.Item(117, 8).Value = 39

$ws.Cells.Item(118, 1).Value = "Suazilandia"
$ws.Cells.Item(118, 2).Value = 2968
$ws.Cells.Item(118, 3).Value = 59
$ws.Cells.Item(118, 4).Value = 1476
$ws.Cells.Item(118, 5).Value = 1437
$ws.Cells.Item(118, 6).Value = 0
$ws.Cells.Item(118, 7).Value = 2
$ws.Cells.Item(118, 8).Value = 55

$ws.Cells.Item(119, 1).Value = "Sri Lanka"
$ws.Cells.Item(119, 2).Value = 2839
$ws.Cells.Item(119, 3).Value = 0
$ws.Cells.Item(119, 4).Value = 2541
$ws.Cells.Item(119, 5).Value = 287
$ws.Cells.Item(119, 6).Value = 0
$ws.Cells.Item(119, 7).Value = 0
$ws.Cells.Item(119, 8).Value = 11

$ws.Cells.Item(122, 1).Value = "Namibia"
$ws.Cells.Item(122, 2).Value = 2652
$ws.Cells.Item(122, 3).Value = 112
$ws.Cells.Item(122, 4).Value = 563
$ws.Cells.Item(122, 5).Value = 2074
$ws.Cells.Item(122, 6).Value = 0
$ws.Cells.Item(122, 7).Value = 3
$ws.Cells.Item(122, 8).Value = 15

$ws.Cells.Item(123, 1).Value = "Mali"
$ws.Cells.Item(123, 2).Value = 2546
$ws.Cells.Item(123, 3).Value = 0
$ws.Cells.Item(123, 4).Value = 1950
$ws.Cells.Item(123, 5).Value = 472
$ws.Cells.Item(123, 6).Value = 0
$ws.Cells.Item(123, 7).Value = 0
$ws.Cells.Item(123, 8).Value = 124

$ws.Cells.Item(125, 1).Value = "Sudan del Sur"
$ws.Cells.Item(125, 2).Value = 2450
$ws.Cells.Item(125, 3).Value = 13
$ws.Cells.Item(125, 4).Value = 1175
$ws.Cells.Item(125, 5).Value = 1228
$ws.Cells.Item(125, 6).Value = 0
$ws.Cells.Item(125, 7).Value = 0
$ws.Cells.Item(125, 8).Value = 47

$ws.Cells.Item(129, 1).Value = "Mozambique"
$ws.Cells.Item(129, 2).Value = 2120
$ws.Cells.Item(129, 3).Value = 41
$ws.Cells.Item(129, 4).Value = 795
$ws.Cells.Item(129, 5).Value = 1310
$ws.Cells.Item(129, 6).Value = 0
$ws.Cells.Item(129, 7).Value = 0
$ws.Cells.Item(129, 8).Value = 15

$ws.Cells.Item(130, 1).Value = "Ruanda"
$ws.Cells.Item(130, 2).Value = 2104
$ws.Cells.Item(130, 3).Value = 0
$ws.Cells.Item(130, 4).Value = 1237
$ws.Cells.Item(130, 5).Value = 862
$ws.Cells.Item(130, 6).Value = 0
$ws.Cells.Item(130, 7).Value = 0
$ws.Cells.Item(130, 8).Value = 5

$ws.Cells.Item(135, 1).Value = "Sierra Leona"
$ws.Cells.Item(135, 2).Value = 1877
$ws.Cells.Item(135, 3).Value = 17
$ws.Cells.Item(135, 4).Value = 1427
$ws.Cells.Item(135, 5).Value = 383
$ws.Cells.Item(135, 6).Value = 0
$ws.Cells.Item(135, 7).Value = 0
$ws.Cells.Item(135, 8).Value = 67

$ws.Cells.Item(150, 1).Value = "Siria"
$ws.Cells.Item(150, 2).Value = 999
$ws.Cells.Item(150, 3).Value = 55
$ws.Cells.Item(150, 4).Value = 311
$ws.Cells.Item(150, 5).Value = 640
$ws.Cells.Item(150, 6).Value = 0
$ws.Cells.Item(150, 7).Value = 0
$ws.Cells.Item(150, 8).Value = 48

$ws.Cells.Item(151, 1).Value = "Malta"
$ws.Cells.Item(151, 2).Value = 946
$ws.Cells.Item(151, 3).Value = 20
$ws.Cells.Item(151, 4).Value = 670
$ws.Cells.Item(151, 5).Value = 267
$ws.Cells.Item(151, 6).Value = 0
$ws.Cells.Item(151, 7).Value = 0
$ws.Cells.Item(151, 8).Value = 9

$ws.Cells.Item(176, 1).Value = "Islas Feroe"
$ws.Cells.Item(176, 2).Value = 266
$ws.Cells.Item(176, 3).Value = 25
$ws.Cells.Item(176, 4).Value = 192
$ws.Cells.Item(176, 5).Value = 74
$ws.Cells.Item(176, 6).Value = 0
$ws.Cells.Item(176, 7).Value = 0
$ws.Cells.Item(176, 8).Value = 0

$ws.Cells.Item(177, 1).Value = "Camboya"
$ws.Cells.Item(177, 2).Value = 243
$ws.Cells.Item(177, 3).Value = 0
$ws.Cells.Item(177, 4).Value = 210
$ws.Cells.Item(177, 5).Value = 33
$ws.Cells.Item(177, 6).Value = 0
$ws.Cells.Item(177, 7).Value = 0
$ws.Cells.Item(177, 8).Value = 0

$ws.Cells.Item(202, 1).Value = "Timor Oriental"
$ws.Cells.Item(202, 2).Value = 25
$ws.Cells.Item(202, 3).Value = 0
$ws.Cells.Item(202, 4).Value = 24
$ws.Cells.Item(202, 5).Value = 1
$ws.Cells.Item(202, 6).Value = 0
$ws.Cells.Item(202, 7).Value = 0
$ws.Cells.Item(202, 8).Value = 0

$ws.Cells.Item(203, 1).Value = "Santa Lucia"
$ws.Cells.Item(203, 2).Value = 25
$ws.Cells.Item(203, 3).Value = 0
$ws.Cells.Item(203, 4).Value = 24
$ws.Cells.Item(203, 5).Value = 1
$ws.Cells.Item(203, 6).Value = 0
$ws.Cells.Item(203, 7).Value = 0
$ws.Cells.Item(203, 8).Value = 0

